# Update DM integration test fixture
#
# 1) Bold the header row of every sheet.
# 2) Widen the "best fit" columns to account for the now-bold header text
#    (bold Arial header cells need ~1.5x the width of the regular-weight
#    header; columns whose width is actually driven by the (unchanged,
#    non-bold) data row below the header are left as-is).
# 3) Refresh the CodeScheme row's ID (GUID) on the CodeSchemes sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CodeSchemes sheet
# ---------------------------------------------------------------------
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")

$wsCodeSchemes.Range("A1:N1").Font.Bold = $true

# Columns A (ID) and G (PREFLABEL_FI) are wider than their bold header
# needs (they are sized for the GUID / description text in row 2), so
# their width is unaffected by the header going bold.
# CODEVALUE: before bold, this column's width (14.3) was driven by the
# data row ("EDA11-2018-1"), which is wider than the regular-weight
# header ("CODEVALUE" needs only 12.1). Once bold, the header needs
# 12.1 * 1.5 = 18.15, which overtakes the (unchanged, non-bold) data
# width, so the column widens to 18.15 rather than 14.3 * 1.5.
$wsCodeSchemes.Columns.Item(2).ColumnWidth = 17.428571428571427   # CODEVALUE
$wsCodeSchemes.Columns.Item(3).ColumnWidth = 25.714285714285715   # INFORMATIONDOMAIN
$wsCodeSchemes.Columns.Item(4).ColumnWidth = 22.428571428571427   # LANGUAGECODE
$wsCodeSchemes.Columns.Item(5).ColumnWidth = 14.142857142857142   # STATUS
$wsCodeSchemes.Columns.Item(6).ColumnWidth = 19.142857142857142   # DEFAULTCODE
$wsCodeSchemes.Columns.Item(8).ColumnWidth = 19.142857142857142   # PREFLABEL_EN
$wsCodeSchemes.Columns.Item(9).ColumnWidth = 20.714285714285715   # DESCRIPTION_FI
$wsCodeSchemes.Columns.Item(10).ColumnWidth = 24.0                # DESCRIPTION_EN
$wsCodeSchemes.Columns.Item(11).ColumnWidth = 19.142857142857142  # STARTDATE
$wsCodeSchemes.Columns.Item(12).ColumnWidth = 15.714285714285714  # ENDDATE
$wsCodeSchemes.Columns.Item(13).ColumnWidth = 20.714285714285715  # CODESSHEET
$wsCodeSchemes.Columns.Item(14).ColumnWidth = 27.285714285714285  # EXTENSIONSSHEET

# Refresh the CodeScheme GUID in row 2
$wsCodeSchemes.Range("A2").Value = "42058b83-d1d7-4018-bd20-ea053d4f9aef"

# ---------------------------------------------------------------------
# Codes sheet
# ---------------------------------------------------------------------
$wsCodes = $wb.Worksheets.Item("Codes")

$wsCodes.Range("A1:J1").Font.Bold = $true

$wsCodes.Columns.Item(1).ColumnWidth = 5.857142857142857    # ID
$wsCodes.Columns.Item(2).ColumnWidth = 17.428571428571427   # CODEVALUE
$wsCodes.Columns.Item(3).ColumnWidth = 15.714285714285714   # BROADER
$wsCodes.Columns.Item(4).ColumnWidth = 14.142857142857142   # STATUS
$wsCodes.Columns.Item(5).ColumnWidth = 15.714285714285714   # PREFLABEL_FI
$wsCodes.Columns.Item(6).ColumnWidth = 19.142857142857142   # PREFLABEL_EN
$wsCodes.Columns.Item(7).ColumnWidth = 20.714285714285715   # DESCRIPTION_FI
$wsCodes.Columns.Item(8).ColumnWidth = 24.0                 # DESCRIPTION_EN
$wsCodes.Columns.Item(9).ColumnWidth = 19.142857142857142   # STARTDATE
$wsCodes.Columns.Item(10).ColumnWidth = 15.714285714285714  # ENDDATE

# ---------------------------------------------------------------------
# Extensions sheet
# ---------------------------------------------------------------------
$wsExtensions = $wb.Worksheets.Item("Extensions")

$wsExtensions.Range("A1:I1").Font.Bold = $true

$wsExtensions.Columns.Item(1).ColumnWidth = 5.857142857142857    # ID
$wsExtensions.Columns.Item(2).ColumnWidth = 17.428571428571427   # CODEVALUE
$wsExtensions.Columns.Item(3).ColumnWidth = 14.142857142857142   # STATUS
$wsExtensions.Columns.Item(4).ColumnWidth = 24.0                 # PROPERTYTYPE
$wsExtensions.Columns.Item(5).ColumnWidth = 15.714285714285714   # PREFLABEL_FI
$wsExtensions.Columns.Item(6).ColumnWidth = 19.142857142857142   # PREFLABEL_EN
$wsExtensions.Columns.Item(7).ColumnWidth = 19.142857142857142   # STARTDATE
$wsExtensions.Columns.Item(8).ColumnWidth = 15.714285714285714   # ENDDATE
$wsExtensions.Columns.Item(9).ColumnWidth = 24.0                 # MEMBERSSHEET
